# Add a new event row to the "Card4" sheet, and backfill any previously blank
# cells in the existing data rows (2-19) with the literal text "nan" (mirrors
# the pandas/openpyxl export behaviour that produced this workbook, where a
# DataFrame's NaN cells are written out as the string "nan").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card4")

$lastDataRow = 19
$lastCol = 15   # column O

for ($r = 2; $r -le $lastDataRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Text -eq "") {
            $cell.Value = "nan"
        }
    }
}

# Append the new event as row 20.
$newRow = $lastDataRow + 1

# Column A holds the card number as text throughout the sheet (e.g. "4"), so
# force a text format before assigning the numeric-looking string to avoid it
# being stored as a number.
$cardCell = $ws.Cells.Item($newRow, 1)
$cardCell.NumberFormat = "@"
$cardCell.Value = "4"

$ws.Cells.Item($newRow, 12).Value = "19/1/2026"
$ws.Cells.Item($newRow, 13).Value = "زياره وكيل"
$ws.Cells.Item($newRow, 14).Value = "تغير سوفت ماكينه لنظام bc"
$ws.Cells.Item($newRow, 15).Value = "م.احمدعلي ترتشلر"
